$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Prescaler value (B6): 4999 -> 9
$ws.Range("B6").Value = 9

# Update Time base Required (B14): 1 -> 0.00001 (10 microseconds)
$ws.Range("B14").Value = 0.00001

# Update the selected cell on the sheet, matching the saved view state
$ws.Range("C18").Select()
